$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 13-16: a new match result rotated in, shifting others down ---
# Row 13
$ws.Range("F13").Value = "Heidenheim"
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = "Hoffenheim"
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 3.33
$ws.Range("K13").Value = "07/08/2023 12:55"
$ws.Range("L13").Value = 3.15
$ws.Range("M13").Value = "26/08/2023 15:28"
$ws.Range("N13").Value = 3.58
$ws.Range("O13").Value = "07/08/2023 12:55"
$ws.Range("P13").Value = 3.83
$ws.Range("Q13").Value = "26/08/2023 15:28"
$ws.Range("R13").Value = 2.1
$ws.Range("S13").Value = "07/08/2023 12:55"
$ws.Range("T13").Value = 2.26
$ws.Range("U13").Value = "26/08/2023 15:29"
$ws.Range("V13").Value = "https://www.betexplorer.com/football/germany/bundesliga/heidenheim-hoffenheim/MXsgFrwj/"

# Row 14
$ws.Range("F14").Value = "Freiburg"
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = "Werder Bremen"
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1.71
$ws.Range("K14").Value = "04/08/2023 16:02"
$ws.Range("L14").Value = 1.77
$ws.Range("M14").Value = "26/08/2023 14:50"
$ws.Range("N14").Value = 3.99
$ws.Range("O14").Value = "04/08/2023 16:02"
$ws.Range("P14").Value = 4.1
$ws.Range("Q14").Value = "26/08/2023 14:50"
$ws.Range("R14").Value = 4.99
$ws.Range("S14").Value = "04/08/2023 16:02"
$ws.Range("T14").Value = 4.67
$ws.Range("U14").Value = "26/08/2023 15:03"
$ws.Range("V14").Value = "https://www.betexplorer.com/football/germany/bundesliga/freiburg-werder-bremen/nBD9N6V9/"

# Row 15
$ws.Range("F15").Value = "Darmstadt"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = "Union Berlin"
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 3.35
$ws.Range("K15").Value = "07/08/2023 12:56"
$ws.Range("L15").Value = 4.17
$ws.Range("M15").Value = "26/08/2023 15:29"
$ws.Range("N15").Value = 3.42
$ws.Range("O15").Value = "07/08/2023 12:56"
$ws.Range("P15").Value = 3.42
$ws.Range("Q15").Value = "26/08/2023 15:23"
$ws.Range("R15").Value = 2.15
$ws.Range("S15").Value = "07/08/2023 12:56"
$ws.Range("T15").Value = 2.04
$ws.Range("U15").Value = "26/08/2023 15:23"
$ws.Range("V15").Value = "https://www.betexplorer.com/football/germany/bundesliga/darmstadt-union-berlin/zyncE2hd/"

# Row 16
$ws.Range("F16").Value = "Bochum"
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = "Dortmund"
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 4.97
$ws.Range("K16").Value = "04/08/2023 16:02"
$ws.Range("L16").Value = 5.62
$ws.Range("M16").Value = "26/08/2023 15:27"
$ws.Range("N16").Value = 4.68
$ws.Range("O16").Value = "04/08/2023 16:02"
$ws.Range("P16").Value = 4.72
$ws.Range("Q16").Value = "26/08/2023 15:29"
$ws.Range("R16").Value = 1.61
$ws.Range("S16").Value = "04/08/2023 16:02"
$ws.Range("T16").Value = 1.57
$ws.Range("U16").Value = "26/08/2023 15:26"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/germany/bundesliga/bochum-dortmund/v9hlGONq/"

# Row 30
$ws.Range("F30").Value = "Wolfsburg"
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = "Union Berlin"
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = 2.19
$ws.Range("K30").Value = "28/08/2023 16:01"
$ws.Range("L30").Value = 2.23
$ws.Range("M30").Value = "16/09/2023 15:27"
$ws.Range("N30").Value = 3.35
$ws.Range("O30").Value = "28/08/2023 16:01"
$ws.Range("P30").Value = 3.35
$ws.Range("Q30").Value = "16/09/2023 15:27"
$ws.Range("R30").Value = 3.33
$ws.Range("S30").Value = "28/08/2023 16:01"
$ws.Range("T30").Value = 3.63
$ws.Range("U30").Value = "16/09/2023 15:29"
$ws.Range("V30").Value = "https://www.betexplorer.com/football/germany/bundesliga/wolfsburg-union-berlin/fcENm3qF/"

# Row 31
$ws.Range("F31").Value = "RB Leipzig"
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = "Augsburg"
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 1.2
$ws.Range("K31").Value = "28/08/2023 16:01"
$ws.Range("L31").Value = 1.27
$ws.Range("M31").Value = "16/09/2023 15:25"
$ws.Range("N31").Value = 7.5
$ws.Range("O31").Value = "28/08/2023 16:01"
$ws.Range("P31").Value = 6.75
$ws.Range("Q31").Value = "16/09/2023 15:19"
$ws.Range("R31").Value = 13.58
$ws.Range("S31").Value = "28/08/2023 16:01"
$ws.Range("T31").Value = 10.25
$ws.Range("U31").Value = "16/09/2023 15:25"
$ws.Range("V31").Value = "https://www.betexplorer.com/football/germany/bundesliga/rb-leipzig-augsburg/l0BFkPE2/"

# Row 32
$ws.Range("F32").Value = "FC Koln"
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = "Hoffenheim"
$ws.Range("I32").Value = 3
$ws.Range("J32").Value = 2.03
$ws.Range("K32").Value = "28/08/2023 16:01"
$ws.Range("L32").Value = 2.22
$ws.Range("M32").Value = "16/09/2023 15:17"
$ws.Range("N32").Value = 3.7
$ws.Range("O32").Value = "28/08/2023 16:01"
$ws.Range("P32").Value = 3.86
$ws.Range("Q32").Value = "16/09/2023 15:27"
$ws.Range("R32").Value = 3.68
$ws.Range("S32").Value = "28/08/2023 16:01"
$ws.Range("T32").Value = 3.19
$ws.Range("U32").Value = "16/09/2023 15:20"
$ws.Range("V32").Value = "https://www.betexplorer.com/football/germany/bundesliga/1-fc-koln-hoffenheim/lbWeVfTm/"

# Row 33
$ws.Range("F33").Value = "Freiburg"
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = "Dortmund"
$ws.Range("I33").Value = 4
$ws.Range("J33").Value = 3.62
$ws.Range("K33").Value = "28/08/2023 16:01"
$ws.Range("L33").Value = 2.97
$ws.Range("M33").Value = "16/09/2023 15:29"
$ws.Range("N33").Value = 3.92
$ws.Range("O33").Value = "28/08/2023 16:01"
$ws.Range("P33").Value = 3.82
$ws.Range("Q33").Value = "16/09/2023 15:27"
$ws.Range("R33").Value = 1.9
$ws.Range("S33").Value = "28/08/2023 16:01"
$ws.Range("T33").Value = 2.36
$ws.Range("U33").Value = "16/09/2023 15:28"
$ws.Range("V33").Value = "https://www.betexplorer.com/football/germany/bundesliga/freiburg-dortmund/0zFJlqU8/"

# Row 34
$ws.Range("F34").Value = "Mainz"
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = "Stuttgart"
$ws.Range("I34").Value = 3
$ws.Range("J34").Value = 2.62
$ws.Range("K34").Value = "28/08/2023 16:01"
$ws.Range("L34").Value = 2.71
$ws.Range("M34").Value = "16/09/2023 15:25"
$ws.Range("N34").Value = 3.4
$ws.Range("O34").Value = "28/08/2023 16:01"
$ws.Range("P34").Value = 3.54
$ws.Range("Q34").Value = "16/09/2023 15:27"
$ws.Range("R34").Value = 2.79
$ws.Range("S34").Value = "28/08/2023 16:01"
$ws.Range("T34").Value = 2.71
$ws.Range("U34").Value = "16/09/2023 15:27"
$ws.Range("V34").Value = "https://www.betexplorer.com/football/germany/bundesliga/mainz-vfb-stuttgart/WMyiWzEs/"

# --- Add new row 38 (copy number formats/styles from row 2 for columns A and E) ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A38").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E38").PasteSpecial(-4122) | Out-Null
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "germany"
$ws.Range("C38").Value = "bundesliga"
$ws.Range("D38").Value = "2023-2024"
$ws.Range("E38").Value = 45191.85416666666
$ws.Range("F38").Value = "Stuttgart"
$ws.Range("G38").Value = 3
$ws.Range("H38").Value = "Darmstadt"
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = 1.48
$ws.Range("K38").Value = "11/09/2023 13:08"
$ws.Range("L38").Value = 1.39
$ws.Range("M38").Value = "22/09/2023 20:29"
$ws.Range("N38").Value = 4.48
$ws.Range("O38").Value = "11/09/2023 13:08"
$ws.Range("P38").Value = 5.25
$ws.Range("Q38").Value = "22/09/2023 20:29"
$ws.Range("R38").Value = 6.09
$ws.Range("S38").Value = "11/09/2023 13:08"
$ws.Range("T38").Value = 8.289999999999999
$ws.Range("U38").Value = "22/09/2023 20:29"
$ws.Range("V38").Value = "https://www.betexplorer.com/football/germany/bundesliga/vfb-stuttgart-darmstadt/bBbotCs6/"
